# Auto-generated edit script: update crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure value cells are treated as plain text (preserve formats like "1.00", "61.839.29")
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45", "B46", "C46", "D46", "E46", "D47", "E47", "B48", "C48", "D48", "E48", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.839.29"
$ws.Range("E2").Value = "  -3.72%  "
$ws.Range("D3").Value = "3.013.90"
$ws.Range("E3").Value = "  -4.41%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "525.03"
$ws.Range("E5").Value = "  -6.52%  "
$ws.Range("D6").Value = "127.08"
$ws.Range("E6").Value = "  -10.31%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.008.25"
$ws.Range("E8").Value = "  -4.26%  "
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "0.147"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").Value = "5.97"
$ws.Range("E11").Value = "  -11.28%  "
$ws.Range("D12").Value = "0.440"
$ws.Range("E12").Value = "  -4.77%  "
$ws.Range("D13").Value = "0.0000217"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "32.80"
$ws.Range("E14").Value = "  -9.69%  "
$ws.Range("D15").Value = "3.512.13"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").Value = "61.804.07"
$ws.Range("E16").Value = "  -3.91%  "
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("D18").Value = "3.008.89"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("D19").Value = "6.40"
$ws.Range("E19").Value = "  -5.71%  "
$ws.Range("D20").Value = "466.71"
$ws.Range("E20").Value = "  -8.33%  "
$ws.Range("D21").Value = "12.86"
$ws.Range("E21").Value = "  -7.78%  "
$ws.Range("D22").Value = "0.675"
$ws.Range("E22").Value = "  -5.17%  "
$ws.Range("D23").Value = "6.80"
$ws.Range("E23").Value = "  -8.13%  "
$ws.Range("D24").Value = "77.62"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").Value = "11.59"
$ws.Range("E25").Value = "  -8.74%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "2.60"
$ws.Range("E27").Value = "  -7.87%  "
$ws.Range("D28").Value = "7.77"
$ws.Range("E28").Value = "  -10.85%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "25.12"
$ws.Range("E30").Value = "  -5.00%  "
$ws.Range("D31").Value = "1.79"
$ws.Range("E31").Value = "  -14.70%  "
$ws.Range("D32").Value = "1.07"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").Value = "55.71"
$ws.Range("E33").Value = "  +3.76%  "
$ws.Range("D34").Value = "2.29"
$ws.Range("E34").Value = "  -11.17%  "
$ws.Range("D35").Value = "5.73"
$ws.Range("E35").Value = "  -5.01%  "
$ws.Range("E36").Value = "  -5.40%  "
$ws.Range("D37").Value = "460.30"
$ws.Range("E37").Value = "  -16.32%  "
$ws.Range("D38").Value = "3.019.60"
$ws.Range("E38").Value = "  -4.46%  "
$ws.Range("D39").Value = "0.0380"
$ws.Range("E39").Value = "  -11.01%  "
$ws.Range("D40").Value = "0.0766"
$ws.Range("E40").Value = "  -6.00%  "
$ws.Range("E41").Value = "  -9.40%  "
$ws.Range("D42").Value = "7.81"
$ws.Range("E42").Value = "  -5.06%  "
$ws.Range("D43").Value = "2.42"
$ws.Range("E43").Value = "  -11.44%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "0.240"
$ws.Range("E45").Value = "  -8.23%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "116.98"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").Value = "0.0₃0511"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "1.94"
$ws.Range("E48").Value = "  -11.04%  "
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("D50").Value = "23.30"
$ws.Range("E50").Value = "  -6.26%  "
$ws.Range("D51").Value = "2.22"
$ws.Range("E51").Value = "  -0.87%  "
